$wb = $excel.ActiveWorkbook

# ==== Sheet: ALC ====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 43334.52
$ws.Range("I62").Value = 65717.875
$ws.Range("K62").Value = 65717.875
$ws.Range("M62").Value = -65093.875
$ws.Range("H65").Value = 43334.52
$ws.Range("I65").Value = 65717.875
$ws.Range("K65").Value = 328589.375
$ws.Range("M65").Value = -325469.375
$ws.Range("H86").Value = 2279.9092
$ws.Range("I86").Value = 2405.0715
$ws.Range("K86").Value = 2405.0715
$ws.Range("M86").Value = -1282.0715
$ws.Range("H89").Value = 2279.9092
$ws.Range("I89").Value = 2405.0715
$ws.Range("K89").Value = 12025.3575
$ws.Range("M89").Value = -6409.3575
$ws.Range("H92").Value = 709.4
$ws.Range("I92").Value = 561.75
$ws.Range("J92").Value = 1300
$ws.Range("K92").Value = 561.75
$ws.Range("L92").Value = 1300
$ws.Range("M92").Value = 686.25
$ws.Range("N92").Value = -3796
$ws.Range("H98").Value = 1092.2667
$ws.Range("J98").Value = 204
$ws.Range("L98").Value = 204
$ws.Range("N98").Value = -3200
$ws.Range("H100").Value = 3859.4
$ws.Range("I100").Value = 3349
$ws.Range("J100").Value = 4625
$ws.Range("K100").Value = 3349
$ws.Range("L100").Value = 4625
$ws.Range("M100").Value = -2808
$ws.Range("N100").Value = -5707
$ws.Range("H103").Value = 45454976
$ws.Range("I103").Value = 390.8
$ws.Range("J103").Value = 83333800
$ws.Range("K103").Value = 1172.4
$ws.Range("L103").Value = 250001400
$ws.Range("M103").Value = -586.4000000000001
$ws.Range("N103").Value = -250002572
$ws.Range("H122").Value = 1092.2667
$ws.Range("J122").Value = 204
$ws.Range("L122").Value = 612
$ws.Range("N122").Value = -5512
$ws.Range("H137").Value = 2225.9333
$ws.Range("I137").Value = 1814.3914
$ws.Range("K137").Value = 5443.174199999999
$ws.Range("M137").Value = -2893.174199999999
$ws.Range("H138").Value = 2609.3667
$ws.Range("I138").Value = 1523.0555
$ws.Range("K138").Value = 4569.166499999999
$ws.Range("M138").Value = 570.8335000000006

# ==== Sheet: ARM ====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2530.8667
$ws.Range("I61").Value = 2229.2
$ws.Range("K61").Value = 2229.2
$ws.Range("M61").Value = -2017.2
$ws.Range("H74").Value = 1133.4445
$ws.Range("I74").Value = 993.5333000000001
$ws.Range("K74").Value = 993.5333000000001
$ws.Range("M74").Value = -119.5333000000001
$ws.Range("H77").Value = 1133.4445
$ws.Range("I77").Value = 993.5333000000001
$ws.Range("K77").Value = 4967.6665
$ws.Range("M77").Value = -599.6665000000003
$ws.Range("H97").Value = 1297.25
$ws.Range("I97").Value = 1297.25
$ws.Range("K97").Value = 1297.25
$ws.Range("M97").Value = -801.25
$ws.Range("H136").Value = 2530.8667
$ws.Range("I136").Value = 2229.2
$ws.Range("K136").Value = 6687.599999999999
$ws.Range("M136").Value = -4137.599999999999

# ==== Sheet: BSM ====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3599.6
$ws.Range("I86").Value = 3249.5
$ws.Range("K86").Value = 3249.5
$ws.Range("M86").Value = -2126.5
$ws.Range("H89").Value = 3599.6
$ws.Range("I89").Value = 3249.5
$ws.Range("K89").Value = 16247.5
$ws.Range("M89").Value = -10631.5
$ws.Range("H94").Value = 863.06665
$ws.Range("J94").Value = 833
$ws.Range("L94").Value = 833
$ws.Range("N94").Value = -1735

# ==== Sheet: CRP ====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1367.8572
$ws.Range("I16").Value = 1519
$ws.Range("J16").Value = 990
$ws.Range("K16").Value = 1519
$ws.Range("L16").Value = 990
$ws.Range("M16").Value = -1232
$ws.Range("N16").Value = -1564
$ws.Range("H88").Value = 14330
$ws.Range("J88").Value = 14330
$ws.Range("L88").Value = 14330
$ws.Range("N88").Value = -15142
$ws.Range("H91").Value = 14330
$ws.Range("J91").Value = 14330
$ws.Range("L91").Value = 14330
$ws.Range("N91").Value = -17138
$ws.Range("H96").Value = 9377.429
$ws.Range("J96").Value = 9377.429
$ws.Range("L96").Value = 9377.429
$ws.Range("N96").Value = -14869.429
$ws.Range("H113").Value = 1367.8572
$ws.Range("I113").Value = 1519
$ws.Range("J113").Value = 990
$ws.Range("K113").Value = 1519
$ws.Range("L113").Value = 990
$ws.Range("M113").Value = 651
$ws.Range("N113").Value = -5330

# ==== Sheet: CUL ====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 57.714287
$ws.Range("I2").Value = 38
$ws.Range("K2").Value = 228
$ws.Range("M2").Value = -115
$ws.Range("H11").Value = 1000
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H117").Value = 4355.3184
$ws.Range("J117").Value = 5298.1177
$ws.Range("L117").Value = 15894.3531
$ws.Range("N117").Value = -22778.3531
$ws.Range("H121").Value = 18520244
$ws.Range("I121").Value = 37037412
$ws.Range("K121").Value = 111112236
$ws.Range("M121").Value = -111110926
$ws.Range("H131").Value = 14895.954
$ws.Range("I131").Value = 683.0625
$ws.Range("J131").Value = 52797
$ws.Range("K131").Value = 2049.1875
$ws.Range("L131").Value = 158391
$ws.Range("M131").Value = 2990.8125
$ws.Range("N131").Value = -168471

# ==== Sheet: GSM ====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3800
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 3800
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H102").Value = 3090.5715
$ws.Range("I102").Value = 2968.1
$ws.Range("K102").Value = 2968.1
$ws.Range("M102").Value = -1346.1
$ws.Range("H113").Value = 3038.1428
$ws.Range("I113").Value = 2636.818
$ws.Range("K113").Value = 2636.818
$ws.Range("M113").Value = -466.8180000000002
$ws.Range("H126").Value = 4763.4736
$ws.Range("I126").Value = 4840.75
$ws.Range("K126").Value = 14522.25
$ws.Range("M126").Value = -12052.25

# ==== Sheet: LTW ====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 13408.926
$ws.Range("I7").Value = 14349.608
$ws.Range("K7").Value = 14349.608
$ws.Range("M7").Value = -14237.608
$ws.Range("H40").Value = 6116.8
$ws.Range("I40").Value = 5647.25
$ws.Range("K40").Value = 5647.25
$ws.Range("M40").Value = -5511.25
$ws.Range("H126").Value = 13408.926
$ws.Range("I126").Value = 14349.608
$ws.Range("K126").Value = 43048.824
$ws.Range("M126").Value = -40578.824

# ==== Sheet: WVR ====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6171.1665
$ws.Range("I81").Value = 6390
$ws.Range("J81").Value = 5077
$ws.Range("K81").Value = 12780
$ws.Range("L81").Value = 10154
$ws.Range("M81").Value = -11719
$ws.Range("N81").Value = -12276
$ws.Range("H84").Value = 6171.1665
$ws.Range("I84").Value = 6390
$ws.Range("J84").Value = 5077
$ws.Range("K84").Value = 63900
$ws.Range("L84").Value = 50770
$ws.Range("M84").Value = -58596
$ws.Range("N84").Value = -61378
$ws.Range("H122").Value = 1839.8223
$ws.Range("I122").Value = 1916.697
$ws.Range("J122").Value = 1628.4166
$ws.Range("K122").Value = 5750.090999999999
$ws.Range("L122").Value = 4885.2498
$ws.Range("M122").Value = -3300.090999999999
$ws.Range("N122").Value = -9785.2498
$ws.Range("H126").Value = 11876
$ws.Range("I126").Value = 6818.2
$ws.Range("J126").Value = 21991.6
$ws.Range("K126").Value = 20454.6
$ws.Range("L126").Value = 65974.79999999999
$ws.Range("M126").Value = -17984.6
$ws.Range("N126").Value = -70914.79999999999
